$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Randomize the proxy:port / proxy_user / proxy_pass per row (previously all 3
# data rows shared the same proxy info — now each row gets its own).
$ws.Range("E2").Value = "102.223.180.59:4000"
$ws.Range("F2").Value = "T9S35vqO2weO"
$ws.Range("G2").Value = "54lUH87CQ968"

$ws.Range("E3").Value = "102.223.180.213:4000"
$ws.Range("F3").Value = "cJ639rTIKg9J"
$ws.Range("G3").Value = "2xIhmaEY21i2"

$ws.Range("F4").Value = "8lw5Ek80aOAU"
$ws.Range("G4").Value = "GsyaAG193XmG"
$ws.Range("E4").Value = "102.223.180.88:4000"

# Widen the data columns so the new, longer values are readable.
$ws.Columns.Item(3).ColumnWidth = 29.666666666666668
$ws.Columns.Item(4).ColumnWidth = 25.666666666666668
$ws.Columns.Item(5).ColumnWidth = 24.5
$ws.Columns.Item(6).ColumnWidth = 23.5
$ws.Columns.Item(7).ColumnWidth = 20.666666666666668

# Move the selection onto the last-edited cell.
$ws.Range("E4").Select()
